# Auto commit at 2025-11-17  8:09:24.30
# Append two new daily rows (2025-11-16) for the two stations to the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (number formats/styles) from the last existing data row (31)
# down onto the two new rows (32 and 33) so the new cells inherit the same
# date / currency / integer formatting used throughout the table.
$ws.Range("A30:F31").Copy()
$ws.Range("A32:F33").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 32: 四方坪站 (Sifangping station)
$ws.Range("A32").Value2 = 45977
$ws.Range("B32").Value = "四方坪站"
$ws.Range("C32").Value2 = 8750.49
$ws.Range("D32").Value2 = 7729.69
$ws.Range("E32").Value2 = 2845.94
$ws.Range("F32").Value2 = 360

# Row 33: 高岭站 (Gaoling station)
$ws.Range("A33").Value2 = 45977
$ws.Range("B33").Value = "高岭站"
$ws.Range("C33").Value2 = 4541.23
$ws.Range("D33").Value2 = 4001.15
$ws.Range("E33").Value2 = 1110
$ws.Range("F33").Value2 = 155

# Match the recorded selection state after the edit.
$ws.Range("I32").Select()
